$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.599.96'
$ws.Range("E2").Value = '  +1.79%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.313.63'
$ws.Range("E3").Value = '  +1.23%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.36'
$ws.Range("E5").Value = '  +0.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.64'
$ws.Range("E6").Value = '  +2.31%  '

$ws.Range("E7").Value = '  +1.00%  '

$ws.Range("E8").Value = '  +0.26%  '

$ws.Range("E9").Value = '  +1.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.09'
$ws.Range("E10").Value = '  +2.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0908'
$ws.Range("E11").Value = '  +0.73%  '

$ws.Range("E12").Value = '  +3.90%  '

$ws.Range("E13").Value = '  +0.90%  '

$ws.Range("E14").Value = '  +2.28%  '

$ws.Range("E15").Value = '  +1.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.662.85'
$ws.Range("E16").Value = '  +1.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.305.17'
$ws.Range("E17").Value = '  +0.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.723.54'
$ws.Range("E18").Value = '  +2.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.61'
$ws.Range("E19").Value = '  +2.29%  '

$ws.Range("E20").Value = '  +1.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.45'
$ws.Range("E21").Value = '  +34.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.90'
$ws.Range("E22").Value = '  +1.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.56'
$ws.Range("E23").Value = '  -1.45%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.38'
$ws.Range("E24").Value = '  -4.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.25'
$ws.Range("E25").Value = '  +0.29%  '

$ws.Range("E26").Value = '  -0.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.89'
$ws.Range("E27").Value = '  +1.35%  '

$ws.Range("E28").Value = '  -2.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.73'
$ws.Range("E29").Value = '  -0.57%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.90'
$ws.Range("E30").Value = '  +5.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '165.67'
$ws.Range("E31").Value = '  +1.73%  '

$ws.Range("E32").Value = '  +7.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0894'
$ws.Range("E33").Value = '  +2.65%  '

$ws.Range("E34").Value = '  -1.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.60'
$ws.Range("E35").Value = '  -8.56%  '

$ws.Range("E36").Value = '  +1.38%  '

$ws.Range("E37").Value = '  +2.36%  '

$ws.Range("E38").Value = '  +1.77%  '

$ws.Range("E39").Value = '  +2.06%  '

$ws.Range("E40").Value = '  -4.47%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.61'
$ws.Range("E41").Value = '  +11.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.62'
$ws.Range("E42").Value = '  -1.74%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.23'
$ws.Range("E43").Value = '  +1.21%  '

$ws.Range("E44").Value = '  +1.23%  '

$ws.Range("E45").Value = '  +0.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.37'
$ws.Range("E46").Value = '  +4.39%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '116.16'
$ws.Range("E47").Value = '  +1.91%  '

$ws.Range("B48").Value = 'ordi'
$ws.Range("C48").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '82.37'
$ws.Range("E48").Value = '  +7.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.90'
$ws.Range("E49").Value = '  -0.29%  '

$ws.Range("E50").Value = '  +0.85%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.624.49'
$ws.Range("E51").Value = '  +5.59%  '
